$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "290.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.29%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.08%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.948"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.16%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07146"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-9.32%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.802"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-13.77%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.687"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.83%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.735"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.66%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8956"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.39%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1642"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.81%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07579"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.41%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08104"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.46%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03062"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.88%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1002"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.38%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001503"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.52%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005751"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.64%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.467"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.02%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.080"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.80%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3278"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.04%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1297"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.39%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.038"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.99%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1998"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.64%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04522"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.54%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.78%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004000"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-10.33%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001248"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.11%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01614"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-7.22%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04361"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-9.00%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007390"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.21%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1306"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.31%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002004"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-15.02%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009119"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-13.05%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006017"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.62%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.20%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "172.72%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002999"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.04%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.20%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.20%"
